$d = $word.ActiveDocument

function Replace-Text($find, $replace, [bool]$matchCase = $true) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $matchCase, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "FAILED to find: $find"
    }
    return $ok
}

# 1) Merge the ": " run with the following quoted-title run (same visible text,
#    Word naturally coalesces same-formatted adjacent runs on replace).
Replace-Text ": „Symulator wyścigów powietrznych w dowolnej scenerii wygenerowanej z mapy wysokościowej terenu.”" `
             ": „Symulator wyścigów powietrznych w dowolnej scenerii wygenerowanej z mapy wysokościowej terenu.”"

# 2) Move the _GoBack bookmark from after "Stachyra Krzysztof" to right after "Dudek Piotr".
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$rName = $d.Content
$rName.Find.Execute("Dudek Piotr", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $rName) | Out-Null

# 3) Merge the split "trzecio" / "-o" / "sobowej..." run into a single run.
Replace-Text "2) Dodanie trzecio-osobowej kamery i niektórych ustawień kamer" `
             "2) Dodanie trzecio-osobowej kamery i niektórych ustawień kamer"

# 4) Piotr Dudek section: insert a new first bullet "Utworzenie pierwszej prezentacji."
$rIns = $d.Content
$rIns.Find.Execute("1) Implementacja w C# w Visual Studio menu startowego z możliwością wyboru z dysku pliku z mapą wysokościową terenu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rIns.Collapse(1) # wdCollapseStart
$rIns.InsertBefore("1) Utworzenie pierwszej prezentacji.`r")

# 5) Renumber "1) Implementacja w C#..." -> "2) ..." and append a trailing period in its own run.
$rNum = $d.Content
$rNum.Find.Execute("1) Implementacja w C# w Visual Studio menu startowego z możliwością wyboru z dysku pliku z mapą wysokościową terenu", $true, $false, $false, $false, $false, $true, 1, $false, "2) Implementacja w C# w Visual Studio menu startowego z możliwością wyboru z dysku pliku z mapą wysokościową terenu", 2) | Out-Null

$rDot = $d.Content
$rDot.Find.Execute("2) Implementacja w C# w Visual Studio menu startowego z możliwością wyboru z dysku pliku z mapą wysokościową terenu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rDot.Collapse(0) # wdCollapseEnd
$rDot.InsertAfter(".")

# 6) Renumber remaining list items (2..8 -> 3..9).
Replace-Text "2) Implementacja odczytu danych wysokościowych" "3) Implementacja odczytu danych wysokościowych"
Replace-Text "3) Implementacja zabezpieczeń" "4) Implementacja zabezpieczeń"
Replace-Text "4) Scalenie utworzonego menu startowego z projektem" "5) Scalenie utworzonego menu startowego z projektem"
Replace-Text "5) Dodanie generacji płaskiego terenu" "6) Dodanie generacji płaskiego terenu"
Replace-Text "6) Dodanie tekstury do generowanego terenu." "7) Dodanie tekstury do generowanego terenu."

# 7) Merge the wrapped two-paragraph items into single paragraphs, with renumbering.
Replace-Text ("7) Dodanie dynamicznego ustawienia rozmiarów terenu na podstawie danych" + [char]13 + "   odczytanych z mapy wysokościowej.") `
             "8) Dodanie dynamicznego ustawienia rozmiarów terenu na podstawie danych odczytanych z mapy wysokościowej."

Replace-Text ("8) Dodanie ustawiania wysokości poszczególnych punktów terenu na podstawie odczytanych z mapy" + [char]13 + "   poziomów szarości (na razie szwankuje).") `
             "9) Dodanie ustawiania wysokości poszczególnych punktów terenu na podstawie odczytanych z mapy poziomów szarości (na razie szwankuje)."

# 8) Append a new final bullet "Utworzenie drugiej prezentacji." at the end of the document.
$rEnd = $d.Content
$rEnd.Collapse(0) # wdCollapseEnd
$rEnd.InsertParagraphAfter()
$rEnd.Collapse(0)
$rEnd.InsertAfter("10) Utworzenie drugiej prezentacji.")

Write-Host "Edit complete"
